# feat: add 2022-Q4 data
#
# Starting layout:  总计 (sheet1), 2022-Q3 (sheet2), 2021-Q4 (sheet3)
# Target layout:    总计 (sheet1), 2022-Q4 (NEW), 2022-Q3, 2021-Q4
#
# Strategy:
#  1. Duplicate the existing "2022-Q3" sheet (Worksheets.Copy) so the new
#     "2022-Q4" sheet is placed right before it and inherits identical
#     formatting/styles, then rename it and grow it from 2 to 5 data rows
#     by copying the styled data-row pattern downward.
#  2. Overwrite all cell values on the new sheet with the 2022-Q4 fund data.
#  3. Update the "总计" (total) sheet: shift the existing two data rows down
#     one row and insert the new 2022-Q4 summary row at the top.
#  4. Re-activate the sheet that was active before (2021-Q4, now the last
#     tab) so the workbook-level active tab / tab selection is preserved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet by copying "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)          # current "2022-Q3"
$q3.Copy($q3)                         # inserts a copy right before $q3
$q4 = $wb.Worksheets.Item(2)          # the freshly inserted copy
$q4.Name = "2022-Q4"

# Grow the sheet from 2 data rows to 5 data rows: replicate the styled
# second row (A2:H2, which carries the s="2" style on column A) down
# into rows 4-6 so every row has matching formatting before we overwrite
# the values.
$styledRow = $q4.Range("A2:H2")
$styledRow.Copy($q4.Range("A4:H4"))
$styledRow.Copy($q4.Range("A5:H5"))
$styledRow.Copy($q4.Range("A6:H6"))

# ---------------------------------------------------------------------
# Step 2: write the 2022-Q4 fund data
# ---------------------------------------------------------------------
# Columns B, D, E, F, G must stay plain text (matching the workbook's
# existing convention of storing these figures as inline/shared strings,
# not numbers) even though they look numeric - force text via "@" number
# format, then drop the visible style back to Normal so no cell ends up
# with a different look than the source data. (Two separate Range calls
# because a single comma-joined multi-area range doesn't reliably apply
# the format to every area.)
$textColB = $q4.Range("B2:B6")
$textColDG = $q4.Range("D2:G6")
$textColB.NumberFormat = "@"
$textColDG.NumberFormat = "@"

$data = @(
    @("0", "004702", "南方金融主题灵活配置混合A",       "12.97", "92.71", "5.10", "0.6615", 4),
    @("1", "013500", "南方金融主题灵活配置混合C",       "4.80",  "92.71", "5.10", "0.2448", 4),
    @("2", "011743", "华夏兴源稳健一年持有期混合A",     "8.78",  "22.90", "0.52", "0.0457", 9),
    @("3", "004321", "前海开源沪港深强国产业灵活配置混合", "0.33",  "90.19", "4.85", "0.0160", 5),
    @("4", "011744", "华夏兴源稳健一年持有期混合C",     "2.75",  "22.90", "0.52", "0.0143", 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $q4.Cells.Item($row, 1).Value = [int]$vals[0]
    $q4.Cells.Item($row, 2).Value = $vals[1]
    $q4.Cells.Item($row, 3).Value = $vals[2]
    $q4.Cells.Item($row, 4).Value = $vals[3]
    $q4.Cells.Item($row, 5).Value = $vals[4]
    $q4.Cells.Item($row, 6).Value = $vals[5]
    $q4.Cells.Item($row, 7).Value = $vals[6]
    $q4.Cells.Item($row, 8).Value = [int]$vals[7]
}

# Drop the temporary "@" number format back to the default look (Normal
# style) now that the text values are safely stored as text.
$textColB.Style = "Normal"
$textColDG.Style = "Normal"

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Row 4 is brand new territory for this sheet (it previously only had
# rows 1-3), so column A needs the same style (s="2") that column A
# already carries on rows 2-3. Clone it from A3 before writing the value.
$total.Range("A3").Copy($total.Range("A4"))

# Push the existing rows 2 and 3 down to rows 3 and 4.
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q4"
$total.Cells.Item(4, 3).Value = 1
$total.Cells.Item(4, 4).Value = 0

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.66

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.98

# ---------------------------------------------------------------------
# Step 4: restore the active tab to what used to be "2021-Q4" (now the
# last sheet), since Copy() above made the new sheet the active one.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

Write-Host "2022-Q4 sheet added; sheets now:" $wb.Worksheets.Count
